$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) updates: force text storage so values like "160.57"
# stay as strings (matching the inline-string cells in the workbook) rather
# than being auto-converted to numbers by Excel, then restore the original
# (default) cell style so no formatting changes are introduced.
$priceUpdates = @{
    "D2" = '65.341.03'
    "D3" = '3.466.99'
    "D5" = '578.13'
    "D6" = '160.57'
    "D7" = '0.612'
    "D8" = '0.999'
    "D9" = '3.466.52'
    "D10" = '7.30'
    "D12" = '0.453'
    "D13" = '4.063.41'
    "D15" = '0.0000192'
    "D16" = '28.44'
    "D17" = '65.264.75'
    "D18" = '3.461.59'
    "D20" = '14.33'
    "D21" = '382.13'
    "D22" = '8.18'
    "D23" = '0.561'
    "D24" = '72.85'
    "D26" = '0.0000121'
    "D27" = '10.13'
    "D28" = '0.179'
    "D29" = '0.998'
    "D31" = '6.20'
    "D32" = '2.06'
    "D33" = '23.65'
    "D36" = '160.76'
    "D37" = '1.93'
    "D38" = '0.0781'
    "D39" = '27.57'
    "D40" = '4.81'
    "D41" = '6.86'
    "D42" = '2.889.08'
    "D44" = '43.29'
    "D45" = '0.784'
    "D46" = '26.18'
    "D47" = '324.07'
    "D50" = '0.880'
    "D51" = '6.65'
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    $cell.Style = "Normal"
}

# Volume(1h) column (E) updates: plain text assignment is sufficient since
# the padded "  +x.xx%  " strings never parse as numbers.
$volumeUpdates = @{
    "E2" = '  +1.38%  '
    "E3" = '  +0.37%  '
    "E4" = '  -0.05%  '
    "E5" = '  +0.53%  '
    "E6" = '  +1.52%  '
    "E7" = '  +5.57%  '
    "E8" = '  -0.07%  '
    "E9" = '  +0.24%  '
    "E10" = '  -1.07%  '
    "E11" = '  +0.44%  '
    "E12" = '  +2.12%  '
    "E13" = '  +0.49%  '
    "E14" = '  +0.50%  '
    "E15" = '  -0.87%  '
    "E16" = '  +0.61%  '
    "E17" = '  +1.29%  '
    "E18" = '  +0.63%  '
    "E19" = '  +0.41%  '
    "E20" = '  -0.57%  '
    "E21" = '  -1.93%  '
    "E22" = '  -1.03%  '
    "E23" = '  +3.18%  '
    "E24" = '  -1.15%  '
    "E25" = '  -0.33%  '
    "E26" = '  +0.04%  '
    "E27" = '  +5.53%  '
    "E28" = '  -0.90%  '
    "E29" = '  -0.83%  '
    "E30" = '  +4.38%  '
    "E31" = '  +0.27%  '
    "E32" = '  +1.17%  '
    "E33" = '  +0.07%  '
    "E34" = '  +4.86%  '
    "E35" = '  +8.20%  '
    "E36" = '  -0.05%  '
    "E37" = '  +2.17%  '
    "E38" = '  +0.97%  '
    "E39" = '  +1.35%  '
    "E40" = '  +8.52%  '
    "E41" = '  +4.45%  '
    "E42" = '  -1.64%  '
    "E43" = '  +0.66%  '
    "E44" = '  +1.45%  '
    "E45" = '  +1.97%  '
    "E46" = '  +10.49%  '
    "E47" = '  +9.64%  '
    "E48" = '  +1.27%  '
    "E49" = '  +2.49%  '
    "E50" = '  +1.84%  '
    "E51" = '  +1.73%  '
}

foreach ($addr in $volumeUpdates.Keys) {
    $ws.Range($addr).Value = $volumeUpdates[$addr]
}

